# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple worksheets, per the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 541.25
$ws.Range("I2").Value = 590.1667
$ws.Range("J2").Value = 394.5
$ws.Range("K2").Value = 590.1667
$ws.Range("L2").Value = 394.5
$ws.Range("M2").Value = -477.1667
$ws.Range("N2").Value = -620.5
$ws.Range("H9").Value = 565750.4
$ws.Range("I9").Value = 735407.9
$ws.Range("K9").Value = 735407.9
$ws.Range("M9").Value = -735238.9
$ws.Range("H29").Value = 894
$ws.Range("J29").Value = 924.3333
$ws.Range("L29").Value = 2772.9999
$ws.Range("N29").Value = -3334.9999
$ws.Range("H39").Value = 689.75
$ws.Range("I39").Value = 865.6
$ws.Range("J39").Value = 396.66666
$ws.Range("K39").Value = 2596.8
$ws.Range("L39").Value = 1189.99998
$ws.Range("M39").Value = -2300.8
$ws.Range("N39").Value = -1781.99998
$ws.Range("H52").Value = 9
$ws.Range("I52").Value = 9
$ws.Range("K52").Value = 27
$ws.Range("M52").Value = 133
$ws.Range("H74").Value = 11986.875
$ws.Range("I74").Value = 13916
$ws.Range("K74").Value = 13916
$ws.Range("M74").Value = -12980
$ws.Range("H77").Value = 11986.875
$ws.Range("I77").Value = 13916
$ws.Range("K77").Value = 69580
$ws.Range("M77").Value = -64900
$ws.Range("H80").Value = 2056.111
$ws.Range("J80").Value = 2267.6667
$ws.Range("L80").Value = 6803.000100000001
$ws.Range("N80").Value = -8799.000100000001
$ws.Range("H83").Value = 2056.111
$ws.Range("J83").Value = 2267.6667
$ws.Range("L83").Value = 20409.0003
$ws.Range("N83").Value = -30393.0003
$ws.Range("H101").Value = 474.75
$ws.Range("I101").Value = 299.66666
$ws.Range("K101").Value = 898.9999799999999
$ws.Range("M101").Value = 723.0000200000001
$ws.Range("H111").Value = 1566.1666
$ws.Range("I111").Value = 849.25
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 2547.75
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = 519.25
$ws.Range("N111").Value = -15134
$ws.Range("H132").Value = 4058.9512
$ws.Range("I132").Value = 4517.5557
$ws.Range("J132").Value = 757
$ws.Range("K132").Value = 13552.6671
$ws.Range("L132").Value = 2271
$ws.Range("M132").Value = -11022.6671
$ws.Range("N132").Value = -7331
$ws.Range("H135").Value = 831.2857
$ws.Range("I135").Value = 674.8889
$ws.Range("J135").Value = 1769.6666
$ws.Range("K135").Value = 6074.0001
$ws.Range("L135").Value = 15926.9994
$ws.Range("M135").Value = -3539.0001
$ws.Range("N135").Value = -20996.9994
$ws.Range("H138").Value = 3739.1943
$ws.Range("I138").Value = 3499.3572
$ws.Range("J138").Value = 3891.818
$ws.Range("K138").Value = 10498.0716
$ws.Range("L138").Value = 11675.454
$ws.Range("M138").Value = -5358.071599999999
$ws.Range("N138").Value = -21955.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3244.7666
$ws.Range("I32").Value = 3244.7666
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3244.7666
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2957.7666
$ws.Range("N32").ClearContents()
$ws.Range("H37").Value = 29533.5
$ws.Range("I37").Value = 29533.5
$ws.Range("K37").Value = 29533.5
$ws.Range("M37").Value = -29260.5
$ws.Range("H45").Value = 20760
$ws.Range("I45").Value = 24906.666
$ws.Range("J45").Value = 2100
$ws.Range("K45").Value = 24906.666
$ws.Range("L45").Value = 2100
$ws.Range("M45").Value = -24529.666
$ws.Range("N45").Value = -2854
$ws.Range("H102").Value = 2006.4286
$ws.Range("I102").Value = 1729.5483
$ws.Range("J102").Value = 4152.25
$ws.Range("K102").Value = 1729.5483
$ws.Range("L102").Value = 4152.25
$ws.Range("M102").Value = -107.5482999999999
$ws.Range("N102").Value = -7396.25
$ws.Range("H122").Value = 6754.077
$ws.Range("I122").Value = 7208.773
$ws.Range("J122").Value = 4253.25
$ws.Range("K122").Value = 21626.319
$ws.Range("L122").Value = 12759.75
$ws.Range("M122").Value = -19176.319
$ws.Range("N122").Value = -17659.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 47622650
$ws.Range("J94").Value = 1606.4615
$ws.Range("L94").Value = 1606.4615
$ws.Range("N94").Value = -2508.4615

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1378.5
$ws.Range("I16").Value = 1462.3077
$ws.Range("J16").Value = 1160.6
$ws.Range("K16").Value = 1462.3077
$ws.Range("L16").Value = 1160.6
$ws.Range("M16").Value = -1175.3077
$ws.Range("N16").Value = -1734.6
$ws.Range("H22").Value = 2499.5
$ws.Range("J22").Value = 2999
$ws.Range("L22").Value = 2999
$ws.Range("N22").Value = -3699
$ws.Range("H113").Value = 1378.5
$ws.Range("I113").Value = 1462.3077
$ws.Range("J113").Value = 1160.6
$ws.Range("K113").Value = 1462.3077
$ws.Range("L113").Value = 1160.6
$ws.Range("M113").Value = 707.6922999999999
$ws.Range("N113").Value = -5500.6
$ws.Range("H134").Value = 3903.2104
$ws.Range("I134").Value = 3333.9565
$ws.Range("K134").Value = 10001.8695
$ws.Range("M134").Value = -7466.869499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 3441.5
$ws.Range("J31").Value = 3441.5
$ws.Range("L31").Value = 10324.5
$ws.Range("N31").Value = -10900.5
$ws.Range("H33").Value = 180.88889
$ws.Range("I33").Value = 176.33333
$ws.Range("J33").Value = 190
$ws.Range("K33").Value = 1057.99998
$ws.Range("L33").Value = 1140
$ws.Range("M33").Value = -774.9999800000001
$ws.Range("N33").Value = -1706
$ws.Range("H44").Value = 2837.1
$ws.Range("I44").Value = 274.4
$ws.Range("J44").Value = 5399.8
$ws.Range("K44").Value = 823.1999999999999
$ws.Range("L44").Value = 16199.4
$ws.Range("M44").Value = -425.1999999999999
$ws.Range("N44").Value = -16995.4
$ws.Range("H140").Value = 10384.5
$ws.Range("J140").Value = 11410.941
$ws.Range("L140").Value = 34232.823
$ws.Range("N140").Value = -44592.823

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4982.778
$ws.Range("I113").Value = 5264.4287
$ws.Range("K113").Value = 5264.4287
$ws.Range("M113").Value = -3094.4287
$ws.Range("H122").Value = 3037.5144
$ws.Range("I122").Value = 2615.5454
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 7846.6362
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -5396.6362
$ws.Range("N122").Value = -34900
$ws.Range("H132").Value = 2141.4138
$ws.Range("I132").Value = 1662.5714
$ws.Range("J132").Value = 3398.375
$ws.Range("K132").Value = 4987.7142
$ws.Range("L132").Value = 10195.125
$ws.Range("M132").Value = -2457.7142
$ws.Range("N132").Value = -15255.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4175.148
$ws.Range("I7").Value = 3588.5386
$ws.Range("J7").Value = 4719.857
$ws.Range("K7").Value = 3588.5386
$ws.Range("L7").Value = 4719.857
$ws.Range("M7").Value = -3476.5386
$ws.Range("N7").Value = -4943.857
$ws.Range("H22").Value = 2887.25
$ws.Range("I22").Value = 3016.6667
$ws.Range("J22").Value = 2499
$ws.Range("K22").Value = 3016.6667
$ws.Range("L22").Value = 2499
$ws.Range("M22").Value = -2721.6667
$ws.Range("N22").Value = -3089
$ws.Range("H27").Value = 2887.25
$ws.Range("I27").Value = 3016.6667
$ws.Range("J27").Value = 2499
$ws.Range("K27").Value = 3016.6667
$ws.Range("L27").Value = 2499
$ws.Range("M27").Value = -2909.6667
$ws.Range("N27").Value = -2713
$ws.Range("H40").Value = 4249.3184
$ws.Range("I40").Value = 3263.8235
$ws.Range("K40").Value = 3263.8235
$ws.Range("M40").Value = -3127.8235
$ws.Range("H46").Value = 1273.75
$ws.Range("I46").Value = 1241.4286
$ws.Range("K46").Value = 1241.4286
$ws.Range("M46").Value = -1053.4286
$ws.Range("H61").Value = 4165.3335
$ws.Range("I61").Value = 4165.3335
$ws.Range("K61").Value = 4165.3335
$ws.Range("M61").Value = -3963.3335
$ws.Range("H68").Value = 6099.5
$ws.Range("I68").Value = 4700
$ws.Range("J68").Value = 7499
$ws.Range("K68").Value = 4700
$ws.Range("L68").Value = 7499
$ws.Range("M68").Value = -3951
$ws.Range("N68").Value = -8997
$ws.Range("H71").Value = 6099.5
$ws.Range("I71").Value = 4700
$ws.Range("J71").Value = 7499
$ws.Range("K71").Value = 23500
$ws.Range("L71").Value = 37495
$ws.Range("M71").Value = -19756
$ws.Range("N71").Value = -44983
$ws.Range("H82").Value = 972.7353000000001
$ws.Range("I82").Value = 951.0345
$ws.Range("J82").Value = 1098.6
$ws.Range("K82").Value = 951.0345
$ws.Range("L82").Value = 1098.6
$ws.Range("M82").Value = -590.0345
$ws.Range("N82").Value = -1820.6
$ws.Range("H85").Value = 972.7353000000001
$ws.Range("I85").Value = 951.0345
$ws.Range("J85").Value = 1098.6
$ws.Range("K85").Value = 951.0345
$ws.Range("L85").Value = 1098.6
$ws.Range("M85").Value = 296.9655
$ws.Range("N85").Value = -3594.6
$ws.Range("H113").Value = 4165.3335
$ws.Range("I113").Value = 4165.3335
$ws.Range("K113").Value = 4165.3335
$ws.Range("M113").Value = -1995.3335
$ws.Range("H126").Value = 4175.148
$ws.Range("I126").Value = 3588.5386
$ws.Range("J126").Value = 4719.857
$ws.Range("K126").Value = 10765.6158
$ws.Range("L126").Value = 14159.571
$ws.Range("M126").Value = -8295.6158
$ws.Range("N126").Value = -19099.571
$ws.Range("H132").Value = 5930
$ws.Range("I132").Value = 2907.4
$ws.Range("K132").Value = 8722.200000000001
$ws.Range("M132").Value = -6192.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2000.1
$ws.Range("I126").Value = 2200.5
$ws.Range("J126").Value = 1699.5
$ws.Range("K126").Value = 6601.5
$ws.Range("L126").Value = 5098.5
$ws.Range("M126").Value = -4131.5
$ws.Range("N126").Value = -10038.5
$ws.Range("H136").Value = 155595.17
$ws.Range("I136").Value = 9176.424000000001
$ws.Range("J136").Value = 615768.4
$ws.Range("K136").Value = 27529.272
$ws.Range("L136").Value = 1847305.2
$ws.Range("M136").Value = -24979.272
$ws.Range("N136").Value = -1852405.2
